$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B63 was stored as text "3" in the original file; convert it to a real number.
$ws.Range("B63").Value = 3

# Append a new row 64 with the additional annotation data.
$ws.Range("A64").Value = "Ying Tang"

# B64 must remain a text value "4" (not a number), so force text formatting
# before assigning it, then drop back to the default style so no stray
# number-format style is left attached to the cell.
$ws.Range("B64").NumberFormat = "@"
$ws.Range("B64").Value = "4"
$ws.Range("B64").Style = "Normal"

$ws.Range("C64").Value = "Take care,"
$ws.Range("D64").Value = "ACK"
$ws.Range("E64").Value = "OTH"
$ws.Range("F64").Value = "a069c61a-6fc8-4b7e-8ba1-d302cd5e6a56"
$ws.Range("G64").Value = "NNP_NfOK_ENK4_annotated.xlsx"
$ws.Range("H64").Value = "Take care,"
